$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.233.16"
$ws.Range("E2").Value = "  +0.50%  "

$ws.Range("D3").Value = "1.859.57"
$ws.Range("E3").Value = "  +0.85%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7029"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.40%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "237.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.31%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07749"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.50%  "

$ws.Range("E9").Value = "  +1.09%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.28"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.44%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08177"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.89%  "

$ws.Range("D12").Value = "1.843.97"
$ws.Range("E12").Value = "  +0.40%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7176"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.59%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.169"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.26%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.16"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.50%  "

$ws.Range("D16").Value = "29.245.18"
$ws.Range("E16").Value = "  +0.09%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.773"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.35%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.37"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.10%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007734"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.39%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "236.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.50%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9987"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.23%  "

$ws.Range("D22").Value = "2.107.63"
$ws.Range("E22").Value = "  -0.10%  "

$ws.Range("E23").Value = "  -0.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.447"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.55%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1473"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.44%  "

$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.997"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.53%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.68"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.00%  "

$ws.Range("E28").Value = "  +0.32%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.012"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.68%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.435"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.51%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.417"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.17%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.480"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.57%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.044"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.25%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05222"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.93%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.169"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.72%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7068"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.20%  "

$ws.Range("E37").Value = "  +0.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.669"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.81%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01843"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.17%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.721"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.80%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9308"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.97%  "

$ws.Range("D42").Value = "1.140.28"
$ws.Range("E42").Value = "  +8.14%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.925"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.74%  "

$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4275"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.24%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "70.78"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.76%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9987"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.12%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "103.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.12%  "

$ws.Range("E48").Value = "  +2.97%  "

$ws.Range("D49").Value = "2.003.71"
$ws.Range("E49").Value = "  +0.15%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.174"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.19%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.979"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.29%  "
